$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric columns (unchanged type, plain numbers) ---
$ws.Range("A2").Value = 102077473
$ws.Range("B2").Value = 96367
$ws.Range("E2").Value = 219874
$ws.Range("Q2").Value = 647720.9098417715
$ws.Range("R2").Value = 6560694.968483768
$ws.Range("S2").Value = 10

# --- Plain text columns (non-numeric-looking strings; safe to set directly) ---
$ws.Range("F2").Value = "Nattviol"
$ws.Range("G2").Value = "Platanthera bifolia"
$ws.Range("H2").Value = "(L.) Rich."
$ws.Range("P2").Value = "Tvetaspåret, Tveta, Srm"
$ws.Range("AW2").Value = "Åsa Johansson"
$ws.Range("AX2").Value = "Åsa Johansson"

# --- Text columns that look numeric / date-like: must be forced to stay text ---
# (set a text NumberFormat first so the engine stores a string instead of
# coercing to a number/date, then clear the format again so no residual
# style index is left on the cell, matching the source file's styling.)
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "3"
$ws.Range("I2").ClearFormats()

$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2022-06-28"
$ws.Range("Y2").ClearFormats()

$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2022-07-05"
$ws.Range("AA2").ClearFormats()

# --- Cell removed entirely in the new version ---
$ws.Range("AI2").Value = ""
